$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11144.286
$ws.Range("J43").Value = 3682.9
$ws.Range("L43").Value = 3682.9
$ws.Range("N43").Value = -3820.9

$ws.Range("I51").Value = 2999.359
$ws.Range("J51").Value = 3128.5789
$ws.Range("K51").Value = 2999.359
$ws.Range("L51").Value = 3128.5789
$ws.Range("M51").Value = -2515.359
$ws.Range("N51").Value = -4096.5789

$ws.Range("H95").Value = 86651.336
$ws.Range("J95").Value = 86651.336
$ws.Range("L95").Value = 86651.336
$ws.Range("N95").Value = -92143.336

$ws.Range("H116").Value = 26869.545
$ws.Range("I116").Value = 4741.25
$ws.Range("K116").Value = 4741.25
$ws.Range("M116").Value = -1299.25

$ws.Range("H135").Value = 993.25
$ws.Range("I135").Value = 667.0952
$ws.Range("J135").Value = 1971.7142
$ws.Range("K135").Value = 6003.8568
$ws.Range("L135").Value = 17745.4278
$ws.Range("M135").Value = -3468.8568
$ws.Range("N135").Value = -22815.4278

$ws.Range("H137").Value = 10631.647
$ws.Range("J137").Value = 20099.25
$ws.Range("L137").Value = 60297.75
$ws.Range("N137").Value = -65397.75

$ws.Range("H138").Value = 2373.8362
$ws.Range("J138").Value = 2527.4314
$ws.Range("L138").Value = 7582.2942
$ws.Range("N138").Value = -17862.2942

$ws.Range("H141").Value = 1972
$ws.Range("I141").Value = 1972
$ws.Range("K141").Value = 5916
$ws.Range("M141").Value = -736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2708.4546
$ws.Range("I45").Value = 1974.125
$ws.Range("K45").Value = 1974.125
$ws.Range("M45").Value = -1597.125

$ws.Range("H63").Value = 2185.611
$ws.Range("I63").Value = 2215.25
$ws.Range("K63").Value = 2215.25
$ws.Range("M63").Value = -1529.25

$ws.Range("H66").Value = 2185.611
$ws.Range("I66").Value = 2215.25
$ws.Range("K66").Value = 11076.25
$ws.Range("M66").Value = -7644.25

$ws.Range("H122").Value = 2886
$ws.Range("I122").Value = 2804.2307
$ws.Range("J122").Value = 3098.6
$ws.Range("K122").Value = 8412.6921
$ws.Range("L122").Value = 9295.799999999999
$ws.Range("M122").Value = -5962.6921
$ws.Range("N122").Value = -14195.8

$ws.Range("H132").Value = 545911.25
$ws.Range("I132").Value = 643451.1
$ws.Range("J132").Value = 2474.7144
$ws.Range("K132").Value = 1930353.3
$ws.Range("L132").Value = 7424.1432
$ws.Range("M132").Value = -1927823.3
$ws.Range("N132").Value = -12484.1432

$ws.Range("H134").Value = 64991
$ws.Range("J134").Value = 64991
$ws.Range("L134").Value = 64991
$ws.Range("N134").Value = -75131

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4564.933
$ws.Range("I134").Value = 2216
$ws.Range("J134").Value = 19833
$ws.Range("K134").Value = 6648
$ws.Range("L134").Value = 59499
$ws.Range("M134").Value = -4113
$ws.Range("N134").Value = -64569

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5403.9165
$ws.Range("I94").Value = 9533.666999999999
$ws.Range("K94").Value = 9533.666999999999
$ws.Range("M94").Value = -9082.666999999999

$ws.Range("H105").Value = 11462.4
$ws.Range("I105").Value = 13215.5
$ws.Range("K105").Value = 13215.5
$ws.Range("M105").Value = -11468.5

$ws.Range("H134").Value = 1992.375
$ws.Range("I134").Value = 2031.8667
$ws.Range("K134").Value = 6095.6001
$ws.Range("M134").Value = -3560.6001

$ws.Range("H141").Value = 813710.3
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 813710.3
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 813710.3
$ws.Range("N141").Value = -824070.3
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3926.5
$ws.Range("J39").Value = 4220.846
$ws.Range("L39").Value = 12662.538
$ws.Range("N39").Value = -13250.538

$ws.Range("H88").Value = 70000
$ws.Range("J88").Value = 40000
$ws.Range("L88").Value = 120000
$ws.Range("N88").Value = -120856

$ws.Range("H91").Value = 70000
$ws.Range("J91").Value = 40000
$ws.Range("L91").Value = 120000
$ws.Range("N91").Value = -122964

$ws.Range("H95").Value = 14000
$ws.Range("J95").Value = 14000
$ws.Range("L95").Value = 42000
$ws.Range("N95").Value = -46118

$ws.Range("H113").Value = 1306.5625
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

$ws.Range("H131").Value = 2903.5625
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 131874.75
$ws.Range("J24").Value = 7856.857
$ws.Range("L24").Value = 7856.857
$ws.Range("N24").Value = -8202.857

$ws.Range("H102").Value = 1778.9354
$ws.Range("I102").Value = 1626.8518
$ws.Range("J102").Value = 2805.5
$ws.Range("K102").Value = 1626.8518
$ws.Range("L102").Value = 2805.5
$ws.Range("M102").Value = -4.851799999999912
$ws.Range("N102").Value = -6049.5

$ws.Range("H122").Value = 4487.1763
$ws.Range("I122").Value = 4290.3335
$ws.Range("K122").Value = 12871.0005
$ws.Range("M122").Value = -10421.0005

$ws.Range("H123").Value = 52963
$ws.Range("J123").Value = 52963
$ws.Range("L123").Value = 52963
$ws.Range("N123").Value = -57863

$ws.Range("H136").Value = 76399.10000000001
$ws.Range("J136").Value = 76399.10000000001
$ws.Range("L136").Value = 229197.3
$ws.Range("N136").Value = -234297.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2433.111
$ws.Range("J22").Value = 2433.111
$ws.Range("L22").Value = 2433.111
$ws.Range("N22").Value = -3023.111

$ws.Range("H27").Value = 2433.111
$ws.Range("J27").Value = 2433.111
$ws.Range("L27").Value = 2433.111
$ws.Range("N27").Value = -2647.111

$ws.Range("H40").Value = 3278.2222
$ws.Range("I40").Value = 1976
$ws.Range("K40").Value = 1976
$ws.Range("M40").Value = -1840

$ws.Range("H55").Value = 1711
$ws.Range("I55").Value = 1620.9231
$ws.Range("K55").Value = 1620.9231
$ws.Range("M55").Value = -1447.9231

$ws.Range("H61").Value = 10687.546
$ws.Range("J61").Value = 15368.167
$ws.Range("L61").Value = 15368.167
$ws.Range("N61").Value = -15772.167

$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751

$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

$ws.Range("H100").Value = 2682
$ws.Range("J100").Value = 3150
$ws.Range("L100").Value = 3150
$ws.Range("N100").Value = -4232

$ws.Range("H113").Value = 10687.546
$ws.Range("J113").Value = 15368.167
$ws.Range("L113").Value = 15368.167
$ws.Range("N113").Value = -19708.167

$ws.Range("H122").Value = 5972
$ws.Range("I122").Value = 4787.0835
$ws.Range("J122").Value = 7749.375
$ws.Range("K122").Value = 14361.2505
$ws.Range("L122").Value = 23248.125
$ws.Range("M122").Value = -11911.2505
$ws.Range("N122").Value = -28148.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 18000
$ws.Range("J51").Value = 18000
$ws.Range("L51").Value = 18000
$ws.Range("N51").Value = -19020

$ws.Range("H52").Value = 11328
$ws.Range("I52").Value = 11328
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 11328
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = -11102
$ws.Range("M52").ClearContents()

$ws.Range("H105").Value = 31500
$ws.Range("J105").Value = 31500
$ws.Range("L105").Value = 31500
$ws.Range("N105").Value = -38488
